# Apply the "adate" prompt-type change and related CRIANCA_VISIT (model) edits.

$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")
$wsPromptTypes = $wb.Worksheets.Item("prompt_types")
$wsModel = $wb.Worksheets.Item("model")

# --- survey sheet: change prompt_type references from custom_date -> adate,
#     and from integer -> text on D18 ---
$wsSurvey.Range("D5").Value = "adate"
$wsSurvey.Range("D18").Value = "text"

# --- prompt_types sheet: define the new "adate" prompt type row ---
$wsPromptTypes.Range("A3").Value = "adate"
$wsPromptTypes.Range("B3").Value = "string"
$wsPromptTypes.Range("C3").Value = "string"
$wsPromptTypes.Range("D3").Value = "Save only mm.dd.yyyy with support for ?? at all positions"

# --- model sheet (CRIANCA_VISIT): update OUTDATE / REGDIA rows to use adate ---
$wsModel.Range("B13").Value = "adate"
$wsModel.Range("B15").Value = "adate"

# --- view/selection state ---
$wsSurvey.Activate()
$wsSurvey.Application.ActiveWindow.ScrollRow = 5
$wsSurvey.Range("D19").Select()

$wb.Save()
